$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table data (columns A..E) for rows 2..57 of the sheet.
# Column order: A (index), B (lamda_1), C (lamda_2), D (dic_nbre_clients_poisson_2_keys),
# E (dic_nbre_clients_prob_poisson_2_values)
$rows = @(
    @(0,33.94444444444444,1.95,0,0.135),
    @(1,33.94444444444444,1.95,2,0.001),
    @(2,33.94444444444444,1.95,3,0.001),
    @(3,33.94444444444444,1.95,4,0.007),
    @(4,33.94444444444444,1.95,5,0.026),
    @(5,33.94444444444444,1.95,6,0.031),
    @(6,33.94444444444444,1.95,7,0.049),
    @(7,33.94444444444444,1.95,8,0.05),
    @(8,33.94444444444444,1.95,9,0.049),
    @(9,33.94444444444444,1.95,10,0.022),
    @(10,33.94444444444444,1.95,11,0.028),
    @(11,33.94444444444444,1.95,12,0.027),
    @(12,33.94444444444444,1.95,13,0.027),
    @(13,33.94444444444444,1.95,14,0.04),
    @(14,33.94444444444444,1.95,15,0.045),
    @(15,33.94444444444444,1.95,16,0.039),
    @(16,33.94444444444444,1.95,17,0.042),
    @(17,33.94444444444444,1.95,18,0.035),
    @(18,33.94444444444444,1.95,19,0.03),
    @(19,33.94444444444444,1.95,20,0.028),
    @(20,33.94444444444444,1.95,21,0.025),
    @(21,33.94444444444444,1.95,22,0.016),
    @(22,33.94444444444444,1.95,23,0.019),
    @(23,33.94444444444444,1.95,24,0.025),
    @(24,33.94444444444444,1.95,25,0.026),
    @(25,33.94444444444444,1.95,26,0.021),
    @(26,33.94444444444444,1.95,27,0.008),
    @(27,33.94444444444444,1.95,28,0.013),
    @(28,33.94444444444444,1.95,29,0.01),
    @(29,33.94444444444444,1.95,30,0.011),
    @(30,33.94444444444444,1.95,31,0.012),
    @(31,33.94444444444444,1.95,32,0.007),
    @(32,33.94444444444444,1.95,33,0.007),
    @(33,33.94444444444444,1.95,34,0.01),
    @(34,33.94444444444444,1.95,35,0.008),
    @(35,33.94444444444444,1.95,36,0.006),
    @(36,33.94444444444444,1.95,37,0.011),
    @(37,33.94444444444444,1.95,38,0.005),
    @(38,33.94444444444444,1.95,39,0.007),
    @(39,33.94444444444444,1.95,40,0.006),
    @(40,33.94444444444444,1.95,41,0.003),
    @(41,33.94444444444444,1.95,42,0.004),
    @(42,33.94444444444444,1.95,43,0.001),
    @(43,33.94444444444444,1.95,44,0.001),
    @(44,33.94444444444444,1.95,45,0.002),
    @(45,33.94444444444444,1.95,46,0.004),
    @(46,33.94444444444444,1.95,47,0.002),
    @(47,33.94444444444444,1.95,48,0.002),
    @(48,33.94444444444444,1.95,49,0.004),
    @(49,33.94444444444444,1.95,50,0.002),
    @(50,33.94444444444444,1.95,51,0.003),
    @(51,33.94444444444444,1.95,52,0.002),
    @(52,33.94444444444444,1.95,53,0.001),
    @(53,33.94444444444444,1.95,55,0.001),
    @(54,33.94444444444444,1.95,56,0.001),
    @(55,33.94444444444444,1.95,62,0.001)
)

# Style template for column A (bold, centered/top-aligned, thin border) lives on A2
# already; new rows beyond the original used range (row 48) need the same format
# copied across so the appended rows match the existing table's look.
$styleTemplate = $ws.Cells.Item(2, 1)

$firstRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $firstRow + $i
    $values = $rows[$i]

    $aCell = $ws.Cells.Item($r, 1)
    if ($r -gt 48) {
        $styleTemplate.Copy() | Out-Null
        $aCell.PasteSpecial(-4122) | Out-Null
    }
    $aCell.Value = $values[0]

    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
}

Write-Output "updated $($rows.Count) rows"
